$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continued bibliography search data extraction: append rows 21-23 (Study_ID 509, 522, 523)

# --- Row 21: Study_ID 509 ---
$ws.Range("A21").Value = 509
$ws.Range("C21").Value = 49
$ws.Range("D21").Value = 'UK'
$ws.Range("E21").Value = 0.5
$ws.Range("F21").Value = 56
$ws.Range("G21").Value = 'UK'
$ws.Range("H21").Value = 'UK'
$ws.Range("I21").Value = 'UK'
$ws.Range("J21").Value = 1
$ws.Range("K21").Value = 'UK'
$ws.Range("L21").Value = 'UK'
$ws.Range("M21").Value = 'UK'
$ws.Range("N21").Value = 'UK'
$ws.Range("O21").Value = 'UK'
$ws.Range("P21").Value = 1
$ws.Range("Q21").Value = 'UK'
$ws.Range("R21").Value = 0.61
$ws.Range("S21").Value = 'UK'
$ws.Range("T21").Value = 'UK'
$ws.Range("U21").Value = 'UK'
$ws.Range("V21").Value = 'UK'
$ws.Range("W21").Value = 'UK'
$ws.Range("X21").Value = 'UK'
$ws.Range("Y21").Value = 'UK'
$ws.Range("Z21").Value = 'UK'
$ws.Range("AA21").Value = 'UK'
$ws.Range("AB21").Value = 'UK'
$ws.Range("AC21").Value = 'UK'
$ws.Range("AD21").Value = 'UK'
$ws.Range("AE21").Value = 'UK'
$ws.Range("AF21").Value = 0.02
$ws.Range("AG21").Value = 'UK'
$ws.Range("AH21").Value = 'UK'
$ws.Range("AI21").Value = 'UK'
$ws.Range("AJ21").Value = 0.449
$ws.Range("AK21").Value = 'UK'
$ws.Range("AL21").Value = 'UK'
$ws.Range("AM21").Value = 'UK'
$ws.Range("AN21").Value = 'UK'
$ws.Range("AO21").Value = 0
$ws.Range("AP21").Value = 0.451
$ws.Range("AQ21").Value = 'UK'
$ws.Range("AR21").Value = 'UK'
$ws.Range("AS21").Value = 'UK'
$ws.Range("AT21").Value = 'UK'
$ws.Range("AU21").Value = 'UK'
$ws.Range("AV21").Value = 'UK'
$ws.Range("AW21").Value = 'UK'
$ws.Range("AX21").Value = 'UK'
$ws.Range("AY21").Value = 'UK'
$ws.Range("AZ21").Value = 'UK'
$ws.Range("BA21").Value = 'UK'
$ws.Range("BB21").Value = 'UK'
$ws.Range("BC21").Value = 0.75
$ws.Range("BD21").Value = 'UK'
$ws.Range("BE21").Value = 'UK'
$ws.Range("BF21").Value = 'UK'
$ws.Range("BG21").Value = 'UK'
$ws.Range("BH21").Value = 'UK'
$ws.Range("BI21").Value = 0.567
$ws.Range("BJ21").Value = 'UK'
$ws.Range("BK21").Value = 'UK'
$ws.Range("BL21").Value = 'UK'
$ws.Range("BM21").Value = 'UK'
$ws.Range("BN21").Value = 'UK'
$ws.Range("BO21").Value = 'UK'
$ws.Range("BP21").Value = 'UK'
$ws.Range("BQ21").Value = 0.531
$ws.Range("BR21").Value = 'UK'
$ws.Range("BS21").Value = 0.388
$ws.Range("BT21").Value = 'UK'
$ws.Range("BU21").Value = 'UK'
$ws.Range("BV21").Value = 'UK'
$ws.Range("BW21").Value = 'UK'
$ws.Range("BX21").Value = 'UK'
$ws.Range("BY21").Value = 'UK'
$ws.Range("BZ21").Value = 'UK'
$ws.Range("CA21").Value = 'UK'
$ws.Range("CB21").Value = 'UK'
$ws.Range("CC21").Value = 'UK'
$ws.Range("CD21").Value = 'UK'
$ws.Range("CE21").Value = 'UK'
$ws.Range("CF21").Value = 0.224
$ws.Range("CG21").Value = 'UK'
$ws.Range("CH21").Value = 0.143

# --- Row 22: Study_ID 522 ---
$ws.Range("A22").Value = 522
$ws.Range("C22").Value = 167
$ws.Range("D22").Value = 167
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 99
$ws.Range("G22").Value = 'UK'
$ws.Range("H22").Value = 'UK'
$ws.Range("I22").Value = 'UK'
$ws.Range("J22").Value = 'UK'
$ws.Range("K22").Value = 106
$ws.Range("L22").Value = 0.637
$ws.Range("M22").Value = 61
$ws.Range("N22").Value = 0.363
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 80
$ws.Range("R22").Value = 0.479
$ws.Range("S22").Value = 70
$ws.Range("T22").Value = 0.419
$ws.Range("U22").Value = 10
$ws.Range("V22").Value = 0.054
$ws.Range("W22").Value = 'UK'
$ws.Range("X22").Value = 'UK'
$ws.Range("Y22").Value = 'UK'
$ws.Range("Z22").Value = 'UK'
$ws.Range("AA22").Value = 10
$ws.Range("AB22").Value = 0.054
$ws.Range("AC22").Value = 'UK'
$ws.Range("AD22").Value = 'UK'
$ws.Range("AE22").Value = 'UK'
$ws.Range("AF22").Value = 'UK'
$ws.Range("AG22").Value = 'UK'
$ws.Range("AH22").Value = 'UK'
$ws.Range("AI22").Value = 'UK'
$ws.Range("AJ22").Value = 'UK'
$ws.Range("AK22").Value = 'UK'
$ws.Range("AL22").Value = 'UK'
$ws.Range("AM22").Value = 'UK'
$ws.Range("AN22").Value = 'UK'
$ws.Range("AO22").Value = 'UK'
$ws.Range("AP22").Value = 1
$ws.Range("AQ22").Value = 'UK'
$ws.Range("AR22").Value = 'UK'
$ws.Range("AS22").Value = 'UK'
$ws.Range("AT22").Value = 'UK'
$ws.Range("AU22").Value = 'UK'
$ws.Range("AV22").Value = 'UK'
$ws.Range("AW22").Value = 167
$ws.Range("AX22").Value = 'UK'
$ws.Range("AY22").Value = 'UK'
$ws.Range("AZ22").Value = 'UK'
$ws.Range("BA22").Value = 'UK'
$ws.Range("BB22").Value = 75
$ws.Range("BC22").Value = 0.449
$ws.Range("BD22").Value = 'UK'
$ws.Range("BE22").Value = 'UK'
$ws.Range("BF22").Value = 103
$ws.Range("BG22").Value = 0.617
$ws.Range("BH22").Value = 'UK'
$ws.Range("BI22").Value = 0.3
$ws.Range("BJ22").Value = 'UK'
$ws.Range("BK22").Value = 'UK'
$ws.Range("BL22").Value = 92
$ws.Range("BM22").Value = 0.551
$ws.Range("BN22").Value = 'UK'
$ws.Range("BO22").Value = 'UK'
$ws.Range("BP22").Value = 100
$ws.Range("BQ22").Value = 0.599
$ws.Range("BR22").Value = 51
$ws.Range("BS22").Value = 0.306
$ws.Range("BT22").Value = 2
$ws.Range("BU22").Value = 0.196
$ws.Range("BV22").Value = 'UK'
$ws.Range("BW22").Value = 'UK'
$ws.Range("BX22").Value = 16
$ws.Range("BY22").Value = 'UK'
$ws.Range("BZ22").Value = 2
$ws.Range("CA22").Value = 'UK'
$ws.Range("CB22").Value = 42
$ws.Range("CC22").Value = 'UK'
$ws.Range("CD22").Value = '"Surgical findings were intraluminal FB in 21 (41.2%) patients, mucosal erosions in 14 (27.5%), and perforation or abscess in 16 (31.3%) patients. Operation was considered to be necessary in 45 (88.2%) patients at the time of surgery. In the 41 patients treated by laparotomy, a gastrotomy was required in 32 patients for retrieval of the ingested FB, a duodenotomy in 1 patient, an enterotomy in 8, and intestinal resection in 2 patients. Complications arose in 9 (17.6%) of the patients undergoing surgery and included wound sepsis in 4, pneumonia in 1, intra-abdominal sepsis in 1, intestinal fistula in 1, and intestinal obstruction in 2 patients. Mortality occurred in 1 (1.9%) patient from aorto-esophageal fistula, a 76 year old woman who had endoscopy showing massive bleeding of an unclear source"'
$ws.Range("CE22").Value = 14
$ws.Range("CF22").Value = 0.084
$ws.Range("CG22").Value = 14
$ws.Range("CH22").Value = 'UK'

# --- Row 23: Study_ID 523 ---
$ws.Range("A23").Value = 523
$ws.Range("B23").Value = '"The most common motive for swallowing was suicidal ideation with command hallucinations, re ported by ten patients. Other pa tients'' motives were recorded as suicidal ideation without command hallucinations (two patients), command hallucinations without sui cidal ideation (two patients), depres sion with a desire to harm but not kill themselves (two patients), and ma nipulation ofthe medicolegal system (three patients).'
$ws.Range("C23").Value = 19
$ws.Range("D23").Value = 19
$ws.Range("E23").Value = 17
$ws.Range("F23").Value = 40
$ws.Range("G23").Value = 24
$ws.Range("H23").Value = 'UK'
$ws.Range("I23").Value = 'UK'
$ws.Range("J23").Value = 1
$ws.Range("K23").Value = 19
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 'UK'
$ws.Range("R23").Value = 1
$ws.Range("S23").Value = 'UK'
$ws.Range("T23").Value = 1
$ws.Range("U23").Value = 0
$ws.Range("V23").Value = 0
$ws.Range("W23").Value = 'UK'
$ws.Range("X23").Value = 'UK'
$ws.Range("Y23").Value = 'UK'
$ws.Range("Z23").Value = 'UK'
$ws.Range("AA23").Value = 18
$ws.Range("AB23").Value = 0.947
$ws.Range("AC23").Value = 'UK'
$ws.Range("AD23").Value = 'UK'
$ws.Range("AE23").Value = 'UK'
$ws.Range("AF23").Value = 'UK'
$ws.Range("AG23").Value = 4
$ws.Range("AH23").Value = 'UK'
$ws.Range("AI23").Value = 3
$ws.Range("AJ23").Value = 'UK'
$ws.Range("AK23").Value = 12
$ws.Range("AL23").Value = 'UK'
$ws.Range("AM23").Value = 0
$ws.Range("AN23").Value = 0
$ws.Range("AO23").Value = 0
$ws.Range("AP23").Value = 0
$ws.Range("AQ23").Value = 0
$ws.Range("AR23").Value = 0
$ws.Range("AS23").Value = 0
$ws.Range("AT23").Value = 0
$ws.Range("AU23").Value = 0
$ws.Range("AV23").Value = 0
$ws.Range("AW23").Value = 'UK'
$ws.Range("AX23").Value = 'UK'
$ws.Range("AY23").Value = 'UK'
$ws.Range("AZ23").Value = 'UK'
$ws.Range("BA23").Value = 'UK'
$ws.Range("BB23").Value = 'UK'
$ws.Range("BC23").Value = 'UK'
$ws.Range("BD23").Value = 'UK'
$ws.Range("BE23").Value = 'KU'
$ws.Range("BF23").Value = 'UK'
$ws.Range("BG23").Value = 1
$ws.Range("BH23").Value = 'UK'
$ws.Range("BI23").Value = 'UK'
$ws.Range("BJ23").Value = 'UK'
$ws.Range("BK23").Value = 'UK'
$ws.Range("BL23").Value = 'UK'
$ws.Range("BM23").Value = 'UK'
$ws.Range("BN23").Value = 'UK'
$ws.Range("BO23").Value = 'UK'
$ws.Range("BP23").Value = 'UK'
$ws.Range("BQ23").Value = 'UK'
$ws.Range("BR23").Value = 5
$ws.Range("BS23").Value = 'UK'
$ws.Range("BT23").Value = 0
$ws.Range("BU23").Value = 'UK'
$ws.Range("BV23").Value = 'UK'
$ws.Range("BW23").Value = 'UK'
$ws.Range("BX23").Value = 'UK'
$ws.Range("BY23").Value = 'UK'
$ws.Range("BZ23").Value = 'UK'
$ws.Range("CA23").Value = 'UK'
$ws.Range("CB23").Value = 'UK'
$ws.Range("CC23").Value = 'UK'
$ws.Range("CD23").Value = 'UK'
$ws.Range("CE23").Value = 14
$ws.Range("CF23").Value = 'UK'
$ws.Range("CG23").Value = 'UK'
$ws.Range("CH23").Value = 'UK'

# Column BF (58) now contains a 3-digit value (row 22); widen it to match the new best-fit width,
# which splits the old single 55-67 run into 55-57 / 58 / 59-67.
$ws.Columns.Item(58).ColumnWidth = 3.33334

# Leave the view scrolled to the top of the frozen pane and select the last cell entered.
$ws.Range("A23").Select()
